# The presentation is already open as $ppt.ActivePresentation.
$p = $ppt.ActivePresentation

# Slide 4 ("Code") has four click-triggered "Appear" animations (one per
# screenshot picture, spid 4-7) that drive a <p:timing> block in the slide
# XML. Remove every animation effect from the slide's timeline so that the
# <p:timing> element is dropped entirely from the saved markup.
$slide = $p.Slides.Item(4)
$timeline = $slide.TimeLine

$mainSeq = $timeline.MainSequence
while ($mainSeq.Count -gt 0) {
    $mainSeq.Item(1).Delete()
}

$interactiveSeqs = $timeline.InteractiveSequences
for ($j = 1; $j -le $interactiveSeqs.Count; $j++) {
    $iseq = $interactiveSeqs.Item($j)
    while ($iseq.Count -gt 0) {
        $iseq.Item(1).Delete()
    }
}
